$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.289.31"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.802.15"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'227.14"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'0.573"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'36.01"
$ws.Range("E8").Value = "  +9.79%  "
$ws.Range("D9").Value = "'0.300"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").Value = "'0.0964"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "2.063.05"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'11.63"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("D14").Value = "1.804.40"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("D15").Value = "'0.643"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").Value = "34.311.00"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'68.93"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'245.02"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "0.0₃0793"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "'11.58"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'4.17"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'172.22"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("E26").Value = "  +8.56%  "
$ws.Range("D27").Value = "'16.82"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'4.01"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'0.0531"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "'3.83"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'1.24"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "1.389.11"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "'0.670"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").Value = "'2.46"
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("D39").Value = "'0.0190"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  +11.02%  "
$ws.Range("D41").Value = "'0.958"
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").Value = "'81.75"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'13.47"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").Value = "'6.02"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("D48").Value = "1.964.59"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "'104.47"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  -0.32%  "
